$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format first so numeric-looking strings (e.g. "590.98")
# are preserved exactly as text like the original inline-string cells, instead of
# being auto-converted to numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "62.849.41"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "3.168.17"
$ws.Range("E3").Value = "  -5.25%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "590.98"
$ws.Range("E5").Value = "  -2.56%  "
$ws.Range("D6").Value = "134.24"
$ws.Range("E6").Value = "  -6.41%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.162.82"
$ws.Range("E8").Value = "  -5.39%  "
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("E10").Value = "  -6.70%  "
$ws.Range("D11").Value = "5.25"
$ws.Range("E11").Value = "  -5.89%  "
$ws.Range("E12").Value = "  -3.94%  "
$ws.Range("D13").Value = "0.0000236"
$ws.Range("E13").Value = "  -5.48%  "
$ws.Range("D14").Value = "34.92"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "3.685.10"
$ws.Range("E15").Value = "  -4.79%  "
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "3.167.85"
$ws.Range("E17").Value = "  -4.71%  "
$ws.Range("D18").Value = "62.858.19"
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("E19").Value = "  -5.13%  "
$ws.Range("D20").Value = "460.05"
$ws.Range("E20").Value = "  -4.91%  "
$ws.Range("D21").Value = "13.89"
$ws.Range("E21").Value = "  -2.23%  "
$ws.Range("D22").Value = "0.697"
$ws.Range("E22").Value = "  -6.05%  "
$ws.Range("D23").Value = "7.64"
$ws.Range("E23").Value = "  -5.10%  "
$ws.Range("D24").Value = "13.40"
$ws.Range("E24").Value = "  -3.96%  "
$ws.Range("D25").Value = "82.99"
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  -4.50%  "
$ws.Range("E29").Value = "  -6.19%  "
$ws.Range("D30").Value = "7.73"
$ws.Range("E30").Value = "  -7.16%  "
$ws.Range("E31").Value = "  -6.49%  "
$ws.Range("D32").Value = "27.13"
$ws.Range("E33").Value = "  -4.96%  "
$ws.Range("D34").Value = "2.36"
$ws.Range("E34").Value = "  -7.47%  "
$ws.Range("E35").Value = "  -6.61%  "
$ws.Range("D36").Value = "5.81"
$ws.Range("E36").Value = "  -5.04%  "
$ws.Range("D37").Value = "51.30"
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("D38").Value = "0.0₃0703"
$ws.Range("E38").Value = "  -7.33%  "
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("D40").Value = "401.11"
$ws.Range("E40").Value = "  -7.48%  "
$ws.Range("D41").Value = "8.09"
$ws.Range("E41").Value = "  -3.52%  "
$ws.Range("E42").Value = "  -4.69%  "
$ws.Range("D43").Value = "2.61"
$ws.Range("E43").Value = "  -6.02%  "
$ws.Range("D44").Value = "2.793.63"
$ws.Range("E44").Value = "  -10.95%  "
$ws.Range("D45").Value = "0.251"
$ws.Range("E45").Value = "  -6.94%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  -6.52%  "
$ws.Range("D48").Value = "124.99"
$ws.Range("E49").Value = "  -5.56%  "
$ws.Range("D50").Value = "34.32"
$ws.Range("E50").Value = "  -5.74%  "
$ws.Range("E51").Value = "  -2.70%  "

# Remove the temporary text formatting so the cell style matches the original
# (no explicit style index), while the stored values remain text.
$dRange.ClearFormats()
